$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.847.72"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.600.05"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.95%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("D9").Value = "2.599.92"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("E13").Value = "  -4.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "3.076.82"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000177"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("D17").Value = "66.813.27"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "2.603.81"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "361.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.07%  "
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "67.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "574.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").Value = "0.0₃0998"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("E31").Value = "  -6.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.08%  "
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -9.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.84%  "
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("E41").Value = "  -5.45%  "
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("E43").Value = "  -4.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.614"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.64%  "

Write-Host "Applied 79 cell updates"
